$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ECB May 2025 update: refresh recalculated forecast values for existing rows
$ws.Cells.Item(5, 2).Value = 3.5908604878022601
$ws.Cells.Item(6, 2).Value = 3.8008905500200849
$ws.Cells.Item(7, 2).Value = 3.6810664065035921
$ws.Cells.Item(9, 2).Value = 3.0374057487693267
$ws.Cells.Item(10, 2).Value = 2.8325894070830571
$ws.Cells.Item(11, 2).Value = 3.4546317652178429
$ws.Cells.Item(13, 2).Value = 3.572021069962886
$ws.Cells.Item(17, 2).Value = 3.5242325153881935
$ws.Cells.Item(18, 3).Value = 2.1500000953674316
$ws.Cells.Item(22, 2).Value = 3.9802982039355528
$ws.Cells.Item(23, 2).Value = 3.8307683257712619
$ws.Cells.Item(24, 2).Value = 5.395415937067412
$ws.Cells.Item(25, 2).Value = 4.17234475307471
$ws.Cells.Item(26, 2).Value = 4.9068829752239349
$ws.Cells.Item(29, 2).Value = 7.5673831320844851
$ws.Cells.Item(30, 2).Value = 6.8514983348259006
$ws.Cells.Item(31, 2).Value = 7.098852398172915
$ws.Cells.Item(32, 2).Value = 6.8909852014929909
$ws.Cells.Item(32, 3).Value = 5.9000000953674316
$ws.Cells.Item(43, 2).Value = 4.185440842636833
$ws.Cells.Item(45, 2).Value = 4.885373206234223
$ws.Cells.Item(45, 3).Value = 4.9000000953674316
$ws.Cells.Item(47, 2).Value = 5.0584100961271883
$ws.Cells.Item(48, 2).Value = 4.5508108926848321
$ws.Cells.Item(49, 2).Value = 4.4845725047123643
$ws.Cells.Item(51, 2).Value = 4.6070537216546565
$ws.Cells.Item(52, 2).Value = 4.3897269615596946
$ws.Cells.Item(53, 2).Value = 4.748305967677882

# Fill in the two previously-empty rows (formatting already present)
$ws.Cells.Item(56, 1).Value = 45474
$ws.Cells.Item(56, 2).Value = 3.9907570458550135
$ws.Cells.Item(56, 3).Value = 3
$ws.Cells.Item(57, 1).Value = 45505
$ws.Cells.Item(57, 2).Value = 4.0688963111175429
$ws.Cells.Item(57, 3).Value = 3

# Append new monthly rows 58-64, copying formatting from row 57 first
$ws.Range("A57:C57").Copy()
$ws.Range("A58:C64").PasteSpecial(-4122)
$excel.CutCopyMode = 0

$ws.Cells.Item(58, 1).Value = 45536
$ws.Cells.Item(58, 2).Value = 3.6999880014156212
$ws.Cells.Item(58, 3).Value = 3
$ws.Cells.Item(59, 1).Value = 45566
$ws.Cells.Item(59, 2).Value = 4.1292384551086423
$ws.Cells.Item(59, 3).Value = 3
$ws.Cells.Item(60, 1).Value = 45597
$ws.Cells.Item(60, 2).Value = 4.5625838084259742
$ws.Cells.Item(60, 3).Value = 3
$ws.Cells.Item(61, 1).Value = 45627
$ws.Cells.Item(61, 2).Value = 4.4608988906167895
$ws.Cells.Item(61, 3).Value = 3.2000000476837158
$ws.Cells.Item(62, 1).Value = 45658
$ws.Cells.Item(62, 2).Value = 4.337500418946691
$ws.Cells.Item(62, 3).Value = 3
$ws.Cells.Item(63, 1).Value = 45689
$ws.Cells.Item(63, 2).Value = 4.8187867992992315
$ws.Cells.Item(63, 3).Value = 4
$ws.Cells.Item(64, 1).Value = 45717
$ws.Cells.Item(64, 2).Value = 4.5453821191276251
$ws.Cells.Item(64, 3).Value = 3.5
